$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Histórico" column (C) entirely - header and all per-row history text
$ws.Range("C1:C6").ClearContents()

# Update existing rows 2-6 with new item names/quantities
$ws.Range("A2").Value = "Marmita"
$ws.Range("B2").Value = 169

$ws.Range("A3").Value = "Porta cantil"
$ws.Range("B3").Value = 28

$ws.Range("A4").Value = "Cantil"
$ws.Range("B4").Value = 19

$ws.Range("A5").Value = "Suspensório"
$ws.Range("B5").Value = 20

$ws.Range("A6").Value = "Cinto"
$ws.Range("B6").Value = 5

# Add new rows 7-12
$ws.Range("A7").Value = "Coldre"
$ws.Range("B7").Value = 5

$ws.Range("A8").Value = "Meia VO"
$ws.Range("B8").Value = 10

$ws.Range("A9").Value = "Meia branca (TFM)"
$ws.Range("B9").Value = 50

$ws.Range("A10").Value = "Saco VO"
$ws.Range("B10").Value = 300

$ws.Range("A11").Value = "Fivela preta"
$ws.Range("B11").Value = 190

$ws.Range("A12").Value = "Fivela dourada"
$ws.Range("B12").Value = 100

# Touch (and reset) row 13/column C formatting so the used range/dimension
# extends to A1:C13 with an empty trailing row, matching the published sheet
$ws.Range("C13").NumberFormat = "General"
$ws.Range("C13").ClearFormats()

# Update the selected cell to I7 (matches the diff's sheetView selection)
$null = $ws.Range("I7").Select()
